$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.293.26'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.043.74'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.63'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('E6').Value = '  -2.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.42'
$ws.Range('E8').Value = '  -3.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.383'
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0782'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.66'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.333.00'
$ws.Range('E13').Value = '  -2.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.59'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.751'
$ws.Range('E15').Value = '  -3.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.27'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.047.69'
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.150.62'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.02'
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.20'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0823'
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.29'
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E25').Value = '  -5.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.64'
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.33'
$ws.Range('E27').Value = '  -3.13%  '
$ws.Range('E28').Value = '  -7.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.92'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('E30').Value = '  -4.43%  '
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('E32').Value = '  -4.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0612'
$ws.Range('E33').Value = '  -3.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.56'
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.23'
$ws.Range('E38').Value = '  -4.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.23'
$ws.Range('E39').Value = '  -3.60%  '
$ws.Range('E40').Value = '  -5.01%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.90'
$ws.Range('E41').Value = '  -0.96%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.474.26'
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.84'
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0939'
$ws.Range('E44').Value = '  -3.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.05'
$ws.Range('E45').Value = '  -5.65%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('E47').Value = '  -3.22%  '
$ws.Range('E48').Value = '  -4.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.09'
$ws.Range('E49').Value = '  -4.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.91'
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.231.15'
$ws.Range('E51').Value = '  -1.71%  '
